$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.248767871177931
$ws.Range("C2").Value = 0.8568015655085476
$ws.Range("D2").Value = 2.374591342905846
$ws.Range("E2").Value = 1.540970909169231
$ws.Range("F2").Value = 1.536197893679839
$ws.Range("G2").Value = 50

$ws.Range("B3").Value = 0.04882215121843843
$ws.Range("C3").Value = 0.7469394130861163
$ws.Range("D3").Value = 1.457578478961404
$ws.Range("E3").Value = 1.207302148992291
$ws.Range("F3").Value = 1.219356202157649
$ws.Range("G3").Value = 47

$ws.Range("B4").Value = 0.05541629173989827
$ws.Range("C4").Value = 0.751460361221017
$ws.Range("D4").Value = 1.308898174481467
$ws.Range("E4").Value = 1.144070878259501
$ws.Range("F4").Value = 1.155355179810744
$ws.Range("G4").Value = 46

$ws.Range("B5").Value = 0.2921665240571426
$ws.Range("C5").Value = 0.8607796199794472
$ws.Range("D5").Value = 2.535407960140988
$ws.Range("E5").Value = 1.592296442293642
$ws.Range("F5").Value = 1.582558739570702
$ws.Range("G5").Value = 46

$ws.Range("B6").Value = 0.07348626165398647
$ws.Range("C6").Value = 0.8093869007750452
$ws.Range("D6").Value = 1.94618247269454
$ws.Range("E6").Value = 1.395056440684225
$ws.Range("F6").Value = 1.409225572003233
$ws.Range("G6").Value = 44

$ws.Range("B7").Value = -0.009028182342375235
$ws.Range("C7").Value = 0.9204870629536859
$ws.Range("D7").Value = 2.016096326641489
$ws.Range("E7").Value = 1.419893068734927
$ws.Range("F7").Value = 1.442583638701003
$ws.Range("G7").Value = 32

$ws.Range("B8").Value = 0.07140803266404031
$ws.Range("C8").Value = 0.918135862785761
$ws.Range("D8").Value = 1.745723223468771
$ws.Range("E8").Value = 1.321258197124533
$ws.Range("F8").Value = 1.341135683000723
$ws.Range("G8").Value = 31

$ws.Range("B9").Value = 0.135685191380836
$ws.Range("C9").Value = 0.7164848718710426
$ws.Range("D9").Value = 0.9226493013901432
$ws.Range("E9").Value = 0.9605463556695967
$ws.Range("F9").Value = 0.9868108548879626
$ws.Range("G9").Value = 14

$ws.Range("B10").Value = 0.3151042569182115
$ws.Range("C10").Value = 0.8922117577591338
$ws.Range("D10").Value = 1.375200986598577
$ws.Range("E10").Value = 1.172689637797903
$ws.Range("F10").Value = 1.207552563243096
$ws.Range("G10").Value = 8
